# Update crypto price/volume data on the active sheet to reflect the
# latest GitHub Actions scrape (Sat Dec 17 21:30:42 UTC 2022).
#
# The "Price" column (D) stores numeric-looking values as TEXT (they were
# written as inline strings, e.g. "236.14", "0.05570", keeping trailing
# zeros / exact formatting). If we just assign a numeric-looking string to
# Range.Value, Excel auto-converts it to a real number (losing trailing
# zeros, e.g. "0.05570" -> 0.0557). To keep these as text we mark the cell
# as Text ("@") before writing the value, then restore the default style
# so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($a1, $value) {
    $rng = $ws.Range($a1)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2"  "236.14"
Set-TextValue "D3"  "21.64"
Set-TextValue "D4"  "5.361"
Set-TextValue "D5"  "0.05570"
Set-TextValue "D7"  "6.455"
Set-TextValue "D8"  "0.7989"
Set-TextValue "D9"  "1.030"
Set-TextValue "D11" "0.07312"
Set-TextValue "D12" "0.03213"
Set-TextValue "D13" "0.02915"
Set-TextValue "D14" "0.09247"
Set-TextValue "D15" "0.001646"
Set-TextValue "D16" "3.258"
Set-TextValue "D17" "0.04761"
Set-TextValue "D18" "0.0005709"
$ws.Range("E18").Value = "17OneONE"
Set-TextValue "D19" "0.006263"
Set-TextValue "D20" "0.005067"
Set-TextValue "D22" "0.0001500"
Set-TextValue "D23" "0.0004181"
Set-TextValue "D24" "3.958"
Set-TextValue "D27" "0.1294"
Set-TextValue "D40" "0.04129"
Set-TextValue "D41" "0.007013"
Set-TextValue "D42" "0.003499"
Set-TextValue "D43" "0.1038"
Set-TextValue "D44" "0.009518"
Set-TextValue "D45" "0.00005440"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "D47" "0.6799"
Set-TextValue "D48" "0.03223"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
Set-TextValue "D49" "0.00002100"
Set-TextValue "D50" "0.01010"
